$d = $word.ActiveDocument

# The document ends with a row describing the "Live Cursor Co-editing"
# feature. Append a new section below it documenting the realtime
# websocket build fix, following the same paragraph pattern used
# throughout the document:
#   (blank) / --- / <title> / Updated: <date> / (blank) / header row / data row
#
# New paragraphs/runs are built with InsertXML so that cell separators
# become real <w:tab/> elements (matching the rest of the document)
# instead of literal tab characters inside the text run.

$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
$rPr = '<w:rPr><w:rFonts w:ascii="Helvetica Light" w:hAnsi="Helvetica Light" w:cs="Helvetica Light"/><w:sz w:val="24"/><w:sz-cs w:val="24"/></w:rPr>'

function Get-CellXml([string]$text) {
    $escaped = $text.Replace("&", "&amp;").Replace("<", "&lt;").Replace(">", "&gt;")
    return '<w:t xml:space="preserve">' + $escaped + '</w:t>'
}

$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range

# --- paragraph 1: blank line ---------------------------------------------
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)

# --- paragraph 2: "---" ---------------------------------------------------
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "---")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)

# --- paragraph 3: "Build Fix Update" --------------------------------------
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "Build Fix Update")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)

# --- paragraph 4: "Updated: 2026-02-18" -----------------------------------
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "Updated: 2026-02-18")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)

# --- paragraph 5: blank line -----------------------------------------------
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)

# --- paragraph 6: header row ------------------------------------------------
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "Module Name") + "<w:tab/>" + (Get-CellXml "Developed") + "<w:tab/>" + (Get-CellXml "Partial Developed") + "<w:tab/>" + (Get-CellXml "Need To Develop")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)

# --- paragraph 7: data row ---------------------------------------------------
$rng = $d.Paragraphs.Item($d.Paragraphs.Count).Range
$rng.InsertParagraphAfter()
$p = $d.Paragraphs.Item($d.Paragraphs.Count)
$p.Range.Collapse(0)
$runInner = $rPr + (Get-CellXml "Realtime Collaboration Build Stability") + "<w:tab/>" + (Get-CellXml "Fixed compile/runtime compatibility for websocket send calls and tenancy warning cleanup") + "<w:tab/>" + (Get-CellXml "-") + "<w:tab/>" + (Get-CellXml "-")
$xml = '<w:p ' + $wNs + '><w:r>' + $runInner + '</w:r></w:p>'
$p.Range.InsertXML($xml)
